# Apply the edits described by the commit:
# "aggiornati i file excel con la durata totale dei singoli test di download e upload"
#
# - Rename the two existing metric headers (B1/C1) to include units (Mb/s)
# - Add two new columns D/E with total download/upload test duration (s)
# - Fill in the three data rows for the new columns
# - Update the chart title (two lines: profile + file size)
# - Resize the data-table columns
# - Reposition/resize the chart
# - Move the active-cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Worksheet data -------------------------------------------------

$ws.Range("B1").Value = "Banda in download (Mb/s)"
$ws.Range("C1").Value = "Banda in upload (Mb/s)"
$ws.Range("D1").Value = "Tempo totale download (s)"
$ws.Range("E1").Value = "Tempo totale upload (s)"

$ws.Range("D2").Value = 19.003
$ws.Range("E2").Value = 16.046

$ws.Range("D3").Value = 16.004
$ws.Range("E3").Value = 19.029

$ws.Range("D4").Value = 16.002
$ws.Range("E4").Value = 16.536

# --- Column widths ----------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 22.5
$ws.Columns.Item(2).ColumnWidth = 22.166666666666668
$ws.Columns.Item(3).ColumnWidth = 23
$ws.Columns.Item(4).ColumnWidth = 22.166666666666668
$ws.Columns.Item(5).ColumnWidth = 21.666666666666668

# --- Chart title --------------------------------------------------------

$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.ChartTitle.Text = "Profilo: 80/10" + "`n" + "Dimensione file: 5MB"

# --- Chart position / size ----------------------------------------------
# Target anchor: from col=A(0) colOff=12700 row=8 rowOff=0
#                to   col=L(11) colOff=622300 row=29 rowOff=63500

$co.Left = 1.0
$co.Top = 128.0
$co.Width = 1024.625
$co.Height = 341.0

# --- Selection ------------------------------------------------------------

$ws.Range("G4").Select() | Out-Null
